$d = $word.ActiveDocument
Write-Host "Fields count:" $d.Fields.Count
Write-Host "Bookmarks count:" $d.Bookmarks.Count
foreach ($b in $d.Bookmarks) {
    Write-Host "bookmark:" $b.Name $b.Range.Start $b.Range.End
}
